$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values are bare percentages need NumberFormat forced to
# Text first, otherwise Excel auto-converts "42%" into the number 0.42.
$percentCells = @("H8", "H20", "H21", "H23", "H26", "H27", "H29", "H34", "H45")
foreach ($pc in $percentCells) {
    $ws.Range($pc).NumberFormat = "@"
}

$ws.Range("E2").Value = "2026-02-27 05:48:37"
$ws.Range("N2").Value = "0.5 °C 5:16 TU"
$ws.Range("E3").Value = "2026-02-27 05:48:40"
$ws.Range("E4").Value = "2026-02-27 05:48:42"
$ws.Range("J4").Value = "1025.9 hPa"
$ws.Range("N4").Value = "6.2 °C 5:01 TU"
$ws.Range("E5").Value = "2026-02-27 05:48:45"
$ws.Range("E6").Value = "2026-02-27 05:48:48"
$ws.Range("J6").Value = "1025.7 hPa"
$ws.Range("N6").Value = "8.9 °C 5:04 TU"
$ws.Range("E7").Value = "2026-02-27 05:48:50"
$ws.Range("E8").Value = "2026-02-27 05:48:53"
$ws.Range("H8").Value = "42%"
$ws.Range("J8").Value = "1025.4 hPa"
$ws.Range("L8").Value = "29.9 km/h - 260º 5:06 TU"
$ws.Range("M8").Value = "12.5 °C 5:27 TU"
$ws.Range("E9").Value = "2026-02-27 05:48:56"
$ws.Range("I9").Value = "0.1 mm"
$ws.Range("M9").Value = "9.1 °C 5:21 TU"
$ws.Range("O9").Value = "8.4 °C"
$ws.Range("E10").Value = "2026-02-27 05:48:58"
$ws.Range("N10").Value = "8.3 °C 5:24 TU"
$ws.Range("O10").Value = "9.2 °C"
$ws.Range("E11").Value = "2026-02-27 05:49:01"
$ws.Range("N11").Value = "1.2 °C 5:22 TU"
$ws.Range("O11").Value = "2.4 °C"
$ws.Range("E12").Value = "2026-02-27 05:49:03"
$ws.Range("M12").Value = "9.8 °C 5:21 TU"
$ws.Range("O12").Value = "8.1 °C"
$ws.Range("E13").Value = "2026-02-27 05:49:06"
$ws.Range("N13").Value = "-3.0 °C 5:08 TU"
$ws.Range("O13").Value = "-1.4 °C"
$ws.Range("E14").Value = "2026-02-27 05:49:08"
$ws.Range("N14").Value = "6.4 °C 5:29 TU"
$ws.Range("O14").Value = "8.8 °C"
$ws.Range("E15").Value = "2026-02-27 05:49:11"
$ws.Range("I15").Value = "0.1 mm"
$ws.Range("O15").Value = "8.4 °C"
$ws.Range("E16").Value = "2026-02-27 05:49:14"
$ws.Range("O16").Value = "2.8 °C"
$ws.Range("E17").Value = "2026-02-27 05:49:16"
$ws.Range("E18").Value = "2026-02-27 05:49:19"
$ws.Range("E19").Value = "2026-02-27 05:49:21"
$ws.Range("O19").Value = "7.9 °C"
$ws.Range("E20").Value = "2026-02-27 05:49:24"
$ws.Range("H20").Value = "60%"
$ws.Range("O20").Value = "1.9 °C"
$ws.Range("E21").Value = "2026-02-27 05:49:27"
$ws.Range("H21").Value = "80%"
$ws.Range("N21").Value = "1.9 °C 5:17 TU"
$ws.Range("O21").Value = "3.7 °C"
$ws.Range("E22").Value = "2026-02-27 05:49:29"
$ws.Range("N22").Value = "-0.4 °C 5:00 TU"
$ws.Range("E23").Value = "2026-02-27 05:49:32"
$ws.Range("H23").Value = "42%"
$ws.Range("E24").Value = "2026-02-27 05:49:35"
$ws.Range("J24").Value = "1026.3 hPa"
$ws.Range("N24").Value = "1.3 °C 5:10 TU"
$ws.Range("O24").Value = "4.3 °C"
$ws.Range("E25").Value = "2026-02-27 05:49:37"
$ws.Range("O25").Value = "4.5 °C"
$ws.Range("E26").Value = "2026-02-27 05:49:40"
$ws.Range("H26").Value = "45%"
$ws.Range("J26").Value = "1024.7 hPa"
$ws.Range("E27").Value = "2026-02-27 05:49:43"
$ws.Range("H27").Value = "46%"
$ws.Range("E28").Value = "2026-02-27 05:49:45"
$ws.Range("N28").Value = "4.4 °C 5:09 TU"
$ws.Range("O28").Value = "5.6 °C"
$ws.Range("E29").Value = "2026-02-27 05:49:48"
$ws.Range("H29").Value = "97%"
$ws.Range("E30").Value = "2026-02-27 05:49:51"
$ws.Range("J30").Value = "1025.6 hPa"
$ws.Range("N30").Value = "9.1 °C 5:29 TU"
$ws.Range("O30").Value = "9.9 °C"
$ws.Range("E31").Value = "2026-02-27 05:49:53"
$ws.Range("E32").Value = "2026-02-27 05:49:56"
$ws.Range("N32").Value = "-1.1 °C 5:29 TU"
$ws.Range("O32").Value = "1.1 °C"
$ws.Range("E33").Value = "2026-02-27 05:49:58"
$ws.Range("N33").Value = "1.0 °C 5:07 TU"
$ws.Range("O33").Value = "2.4 °C"
$ws.Range("E34").Value = "2026-02-27 05:50:01"
$ws.Range("H34").Value = "46%"
$ws.Range("M34").Value = "4.3 °C 5:27 TU"
$ws.Range("O34").Value = "2.0 °C"
$ws.Range("E35").Value = "2026-02-27 05:50:04"
$ws.Range("J35").Value = "1025.4 hPa"
$ws.Range("N35").Value = "8.7 °C 5:18 TU"
$ws.Range("E36").Value = "2026-02-27 05:50:06"
$ws.Range("O36").Value = "9.2 °C"
$ws.Range("E37").Value = "2026-02-27 05:50:09"
$ws.Range("N37").Value = "1.8 °C 5:28 TU"
$ws.Range("O37").Value = "2.7 °C"
$ws.Range("E38").Value = "2026-02-27 05:50:12"
$ws.Range("N38").Value = "6.6 °C 5:28 TU"
$ws.Range("O38").Value = "7.6 °C"
$ws.Range("E39").Value = "2026-02-27 05:50:14"
$ws.Range("E40").Value = "2026-02-27 05:50:17"
$ws.Range("O40").Value = "2.0 °C"
$ws.Range("E41").Value = "2026-02-27 05:50:19"
$ws.Range("J41").Value = "1025.9 hPa"
$ws.Range("O41").Value = "8.7 °C"
$ws.Range("E42").Value = "2026-02-27 05:50:22"
$ws.Range("O42").Value = "8.4 °C"
$ws.Range("E43").Value = "2026-02-27 05:50:24"
$ws.Range("O43").Value = "4.0 °C"
$ws.Range("E44").Value = "2026-02-27 05:50:27"
$ws.Range("O44").Value = "-0.4 °C"
$ws.Range("E45").Value = "2026-02-27 05:50:30"
$ws.Range("H45").Value = "52%"
$ws.Range("J45").Value = "1026.5 hPa"
$ws.Range("L45").Value = "20.5 km/h - 133º 5:23 TU"
$ws.Range("E46").Value = "2026-02-27 05:50:32"
$ws.Range("O46").Value = "6.7 °C"
